# Change WDCP: add FakeAP request and response
# - Insert a new row into the "request_type" table on Sheet1 (A14:C16)
#   for REQ_TYPE_FAKE_AP / 0x02 / "获取当前存在FakeAP威胁的AP的信息".
# - This pushes the following "encrypt_type" table down by one row.
# - Update the active selection to E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 18 (shifts the encrypt_type table, which
# used to start at row 18, down to start at row 19; old rows 18-21 -> 19-22).
$ws.Rows.Item(18).Insert()

# Copy the formatting from the row above (REQ_TYPE_AP_LIST) down onto the
# freshly inserted row so the new row matches the rest of the table.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the newly inserted row 17 with the new request_type entry.
$ws.Range("A17").Value = "REQ_TYPE_FAKE_AP"
$ws.Range("B17").Value = "0x02"
$ws.Range("C17").Value = "获取当前存在FakeAP威胁的AP的信息"

# Update the selection to match the post-edit state.
$ws.Range("E11").Select()
